$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 505
$ws.Range("C2").Value = 252.5556611777037
$ws.Range("D2").Value = 1000
$ws.Range("E2").Value = 10

$ws.Range("B3").Value = 0.5004999999999999
$ws.Range("C3").Value = 0.2548516217338647
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0.001

$ws.Range("B4").Value = 0.001415
$ws.Range("C4").Value = 0.0002984748723009226
$ws.Range("D4").Value = 0.002
$ws.Range("E4").Value = 0.00083
